# Auto-generated edit script applying numeric updates to market-price derived
# columns (H-N) across all 8 leve-profit sheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N113").Value = -9158
$ws.Range("J113").Value = 2650
$ws.Range("K113").Value = 3112.5
$ws.Range("I113").Value = 3112.5
$ws.Range("M113").Value = 141.5
$ws.Range("L113").Value = 2650
$ws.Range("H113").Value = 2958.3333
$ws.Range("L123").Value = 56094.285
$ws.Range("J123").Value = 56094.285
$ws.Range("H123").Value = 56094.285
$ws.Range("N123").Value = -65894.285
$ws.Range("L126").Value = 43432.5
$ws.Range("N126").Value = -53312.5
$ws.Range("H126").Value = 43432.5
$ws.Range("J126").Value = 43432.5
$ws.Range("H127").Value = 787.4286
$ws.Range("M127").Value = 4116.571449999999
$ws.Range("I127").Value = 281.14285
$ws.Range("L127").Value = 5400
$ws.Range("N127").Value = -15320
$ws.Range("J127").Value = 1800
$ws.Range("K127").Value = 843.4285500000001
$ws.Range("J128").Value = 49560
$ws.Range("H128").Value = 49560
$ws.Range("L128").Value = 49560
$ws.Range("N128").Value = -59520
$ws.Range("I138").Value = 1873.0834
$ws.Range("M138").Value = -479.2502000000004
$ws.Range("H138").Value = 3319.5862
$ws.Range("N138").Value = -21891.905
$ws.Range("L138").Value = 11611.905
$ws.Range("J138").Value = 3870.635
$ws.Range("K138").Value = 5619.2502

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 3961.9148
$ws.Range("N32").Value = -8744
$ws.Range("I32").Value = 3961.9148
$ws.Range("J32").Value = 8170
$ws.Range("M32").Value = -3674.9148
$ws.Range("L32").Value = 8170
$ws.Range("H32").Value = 4214.4
$ws.Range("H92").Value = 32029.8
$ws.Range("J92").Value = 32029.8
$ws.Range("L92").Value = 32029.8
$ws.Range("N92").Value = -37021.8
$ws.Range("N122").Value = -12370
$ws.Range("J122").Value = 2490
$ws.Range("I122").Value = 19310.182
$ws.Range("L122").Value = 7470
$ws.Range("M122").Value = -55480.546
$ws.Range("H122").Value = 16722.46
$ws.Range("K122").Value = 57930.546
$ws.Range("H127").Value = 32835.555
$ws.Range("L127").Value = 32835.555
$ws.Range("N127").Value = -42755.555
$ws.Range("J127").Value = 32835.555
$ws.Range("L129").Value = 42889.5
$ws.Range("N129").Value = -52889.5
$ws.Range("H129").Value = 42889.5
$ws.Range("J129").Value = 42889.5
$ws.Range("H132").Value = 1811.7273
$ws.Range("M132").Value = -2048.1155
$ws.Range("J132").Value = 2872.8572
$ws.Range("K132").Value = 4578.1155
$ws.Range("L132").Value = 8618.571599999999
$ws.Range("I132").Value = 1526.0385
$ws.Range("N132").Value = -13678.5716
$ws.Range("N134").Value = -50836
$ws.Range("L134").Value = 40696
$ws.Range("J134").Value = 40696
$ws.Range("H134").Value = 40696
$ws.Range("J135").Value = 41967.11
$ws.Range("H135").Value = 41967.11
$ws.Range("N135").Value = -52107.11
$ws.Range("L135").Value = 41967.11

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N16").Value = -20306.666
$ws.Range("I16").Value = 721.6
$ws.Range("M16").Value = -551.6
$ws.Range("H16").Value = 11218.909
$ws.Range("K16").Value = 721.6
$ws.Range("J16").Value = 19966.666
$ws.Range("L16").Value = 19966.666
$ws.Range("L22").Value = 875
$ws.Range("N22").Value = -1221
$ws.Range("I22").Value = 397.5
$ws.Range("J22").Value = 875
$ws.Range("M22").Value = -224.5
$ws.Range("H22").Value = 636.25
$ws.Range("K22").Value = 397.5
$ws.Range("N122").Value = -52260
$ws.Range("J122").Value = 42460
$ws.Range("L122").Value = 42460
$ws.Range("H122").Value = 42460
$ws.Range("H124").Value = 34570
$ws.Range("N124").Value = -44390
$ws.Range("J124").Value = 34570
$ws.Range("L124").Value = 34570
$ws.Range("L126").Value = 30945.715
$ws.Range("N126").Value = -40825.715
$ws.Range("H126").Value = 30945.715
$ws.Range("J126").Value = 30945.715
$ws.Range("M134").Value = -882.7019999999998
$ws.Range("N134").Value = -17283.9999
$ws.Range("K134").Value = 3417.702
$ws.Range("L134").Value = 12213.9999
$ws.Range("J134").Value = 4071.3333
$ws.Range("H134").Value = 1471.1698
$ws.Range("I134").Value = 1139.234

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N16").Value = -3809.3333
$ws.Range("I16").Value = 9814425
$ws.Range("M16").Value = -9814138
$ws.Range("H16").Value = 4908830
$ws.Range("K16").Value = 9814425
$ws.Range("J16").Value = 3235.3333
$ws.Range("L16").Value = 3235.3333
$ws.Range("M56").ClearContents()
$ws.Range("J56").Value = 8000
$ws.Range("H56").Value = 8000
$ws.Range("K56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("N56").Value = -9690
$ws.Range("L56").Value = 8000
$ws.Range("N113").Value = -7575.3333
$ws.Range("J113").Value = 3235.3333
$ws.Range("K113").Value = 9814425
$ws.Range("I113").Value = 9814425
$ws.Range("M113").Value = -9812255
$ws.Range("L113").Value = 3235.3333
$ws.Range("H113").Value = 4908830

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1687.7778
$ws.Range("M92").Value = -2832
$ws.Range("I92").Value = 1360
$ws.Range("J92").Value = 1753.3334
$ws.Range("L92").Value = 5260.0002
$ws.Range("K92").Value = 4080
$ws.Range("N92").Value = -7756.0002
$ws.Range("N113").Value = -2599330.7
$ws.Range("J113").Value = 864996.9
$ws.Range("K113").Value = 1778.1429
$ws.Range("I113").Value = 592.7143
$ws.Range("M113").Value = 391.8571000000002
$ws.Range("L113").Value = 2594990.7
$ws.Range("H113").Value = 416787.28
$ws.Range("J131").Value = 970.2619
$ws.Range("N131").Value = -12990.7857
$ws.Range("L131").Value = 2910.7857
$ws.Range("H131").Value = 908.59375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 35785
$ws.Range("N62").Value = -37157
$ws.Range("J62").Value = 35785
$ws.Range("L62").Value = 35785
$ws.Range("L65").Value = 107355
$ws.Range("N65").Value = -114219
$ws.Range("H65").Value = 35785
$ws.Range("J65").Value = 35785
$ws.Range("K102").Value = 2090
$ws.Range("H102").Value = 2741.6667
$ws.Range("J102").Value = 6000
$ws.Range("L102").Value = 6000
$ws.Range("I102").Value = 2090
$ws.Range("M102").Value = -468
$ws.Range("N102").Value = -9244
$ws.Range("N113").Value = -6596.2
$ws.Range("J113").Value = 2256.2
$ws.Range("K113").Value = 1917.9286
$ws.Range("I113").Value = 1917.9286
$ws.Range("M113").Value = 252.0714
$ws.Range("L113").Value = 2256.2
$ws.Range("H113").Value = 2006.9474
$ws.Range("N130").Value = -47534.285
$ws.Range("L130").Value = 37494.285
$ws.Range("H130").Value = 37494.285
$ws.Range("J130").Value = 37494.285
$ws.Range("J133").Value = 33670
$ws.Range("N133").Value = -43790
$ws.Range("L133").Value = 33670
$ws.Range("H133").Value = 33670

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N40").Value = -3772
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -1895.6666
$ws.Range("H40").Value = 2241.4285
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 2031.6666
$ws.Range("I40").Value = 2031.6666
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 3400
$ws.Range("K82").Value = 2000
$ws.Range("M82").Value = -1639
$ws.Range("L82").Value = 3400
$ws.Range("H82").Value = 3088.889
$ws.Range("N82").Value = -4122
$ws.Range("H85").Value = 3088.889
$ws.Range("K85").Value = 2000
$ws.Range("J85").Value = 3400
$ws.Range("I85").Value = 2000
$ws.Range("N85").Value = -5896
$ws.Range("L85").Value = 3400
$ws.Range("M85").Value = -752
$ws.Range("L93").Value = 2300
$ws.Range("I93").Value = 1481.4445
$ws.Range("K93").Value = 1481.4445
$ws.Range("N93").Value = -4796
$ws.Range("H93").Value = 1630.2727
$ws.Range("J93").Value = 2300
$ws.Range("M93").Value = -233.4445000000001
$ws.Range("N122").Value = -21873.571
$ws.Range("J122").Value = 5657.857
$ws.Range("I122").Value = 4472
$ws.Range("L122").Value = 16973.571
$ws.Range("M122").Value = -10966
$ws.Range("H122").Value = 5064.9287
$ws.Range("K122").Value = 13416
$ws.Range("N129").Value = -52000
$ws.Range("L129").Value = 42000
$ws.Range("H129").Value = 42000
$ws.Range("J129").Value = 42000
$ws.Range("J133").Value = 82800
$ws.Range("N133").Value = -87860
$ws.Range("L133").Value = 82800
$ws.Range("H133").Value = 82800

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9007.375
$ws.Range("N62").Value = -10906.429
$ws.Range("J62").Value = 9658.429
$ws.Range("L62").Value = 9658.429
$ws.Range("N64").Value = -42096
$ws.Range("H64").Value = 41600
$ws.Range("J64").Value = 41600
$ws.Range("L64").Value = 41600
$ws.Range("L65").Value = 48292.145
$ws.Range("J65").Value = 9658.429
$ws.Range("N65").Value = -54532.145
$ws.Range("H65").Value = 9007.375
$ws.Range("N67").Value = -43316
$ws.Range("H67").Value = 41600
$ws.Range("L67").Value = 41600
$ws.Range("J67").Value = 41600
$ws.Range("N105").Value = -50740.5
$ws.Range("J105").Value = 43752.5
$ws.Range("L105").Value = 43752.5
$ws.Range("H105").Value = 43752.5
$ws.Range("J128").Value = 351810
$ws.Range("H128").Value = 351810
$ws.Range("L128").Value = 351810
$ws.Range("N128").Value = -361770
$ws.Range("N130").Value = -35183
$ws.Range("L130").Value = 25143
$ws.Range("H130").Value = 25143
$ws.Range("J130").Value = 25143
$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 3028.5908
$ws.Range("M132").Value = -6555.7724
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9085.7724
$ws.Range("L132").Value = 0
$ws.Range("I132").Value = 3028.5908

